$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.238.55"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.614.15"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D6").Value = "302.82"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "0.3778"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "51.69"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "0.3528"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("D10").Value = "0.08101"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "22.18"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("D14").Value = "6.365"
$ws.Range("E14").Value = "  -2.70%  "
$ws.Range("D15").Value = "7.267"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").Value = "1.599.91"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "93.96"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "0.06905"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "6.483"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "17.20"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "12.33"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D24").Value = "23.236.14"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "2.513"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("D26").Value = "3.004"
$ws.Range("E26").Value = "  -6.85%  "
$ws.Range("D27").Value = "20.87"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value = "151.12"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").Value = "5.243"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "132.09"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "1.772.43"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").Value = "1.062"
$ws.Range("E32").Value = "  +11.59%  "
$ws.Range("D33").Value = "6.460"
$ws.Range("E33").Value = "  -4.79%  "
$ws.Range("D34").Value = "2.092"
$ws.Range("E34").Value = "  -9.16%  "
$ws.Range("D35").Value = "11.37"
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("D37").Value = "0.08681"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").Value = "0.2453"
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").Value = "0.06931"
$ws.Range("E39").Value = "  -3.65%  "
$ws.Range("D40").Value = "5.843"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D41").Value = "1.326"
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").Value = "0.6879"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("D44").Value = "15.18"
$ws.Range("E44").Value = "  -6.35%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "0.6316"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").Value = "2.249"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("D49").Value = "0.07865"
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("D50").Value = "127.37"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").Value = "1.170"
$ws.Range("E51").Value = "  -2.77%  "
